$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as plain text (e.g. "26.650.50").
# When the new value still looks like a plain number (e.g. "211.40"),
# force Text format first so Excel does not silently convert it to a
# Number (which would drop the trailing zero / change the cell type).

$ws.Range("D2").Value = '26.650.82'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.596.57'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.40'
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("E6").Value = '  +0.64%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.52'
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0838'
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").Value = '1.821.20'
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = '1.608.02'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '26.618.83'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = '0.0₃0736'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.90'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("E21").Value = '  +3.73%  '
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.93'
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.27'
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  +1.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").Value = '1.281.76'
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.619'
$ws.Range("E35").Value = '  -7.05%  '
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("E38").Value = '  -0.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.51'
$ws.Range("E41").Value = '  +2.44%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.784'
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.74'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '1.733.84'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.62'
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("E47").Value = '  -3.27%  '
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.37'
$ws.Range("E51").Value = '  -2.01%  '

# Rows 39/40 got reordered by rank: ARBITRUM moved up to row 39,
# WEMIXToken moved down to row 40 (coin/link/price/volume all swap).
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.834'
$ws.Range("E39").Value = '  -1.11%  '

$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.06'
$ws.Range("E40").Value = '  +20.48%  '
